# Apply updated Leve pricing/profit figures per sheet (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 957.3077
$ws.Range("I9").Value = 1006.125
$ws.Range("J9").Value = 879.2
$ws.Range("K9").Value = 1006.125
$ws.Range("L9").Value = 879.2
$ws.Range("M9").Value = -837.125
$ws.Range("N9").Value = -1217.2
$ws.Range("H17").Value = 1566.8032
$ws.Range("J17").Value = 1566.8032
$ws.Range("L17").Value = 4700.4096
$ws.Range("N17").Value = -5036.4096
$ws.Range("H28").Value = 1096.0667
$ws.Range("I28").Value = 880.1667
$ws.Range("K28").Value = 880.1667
$ws.Range("M28").Value = -395.1667
$ws.Range("H105").Value = 54999
$ws.Range("J105").Value = 54999
$ws.Range("L105").Value = 54999
$ws.Range("N105").Value = -61987
$ws.Range("H131").Value = 1462.25
$ws.Range("I131").Value = 1462.25
$ws.Range("K131").Value = 4386.75
$ws.Range("M131").Value = 653.25
$ws.Range("H137").Value = 1180.909
$ws.Range("I137").Value = 1249.6666
$ws.Range("J137").Value = 871.5
$ws.Range("K137").Value = 3748.9998
$ws.Range("L137").Value = 2614.5
$ws.Range("M137").Value = -1198.9998
$ws.Range("N137").Value = -7714.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2984.3547
$ws.Range("I32").Value = 1467.2963
$ws.Range("K32").Value = 1467.2963
$ws.Range("M32").Value = -1180.2963
$ws.Range("H61").Value = 1353.6957
$ws.Range("I61").Value = 1292.1428
$ws.Range("K61").Value = 1292.1428
$ws.Range("M61").Value = -1080.1428
$ws.Range("H80").Value = 67852.14
$ws.Range("J80").Value = 67852.14
$ws.Range("L80").Value = 67852.14
$ws.Range("N80").Value = -69848.14
$ws.Range("H83").Value = 67852.14
$ws.Range("J83").Value = 67852.14
$ws.Range("L83").Value = 203556.42
$ws.Range("N83").Value = -213540.42
$ws.Range("H122").Value = 2692.5
$ws.Range("I122").Value = 2692.5
$ws.Range("K122").Value = 8077.5
$ws.Range("M122").Value = -5627.5
$ws.Range("H125").Value = 79490.836
$ws.Range("J125").Value = 79490.836
$ws.Range("L125").Value = 79490.836
$ws.Range("N125").Value = -89330.836
$ws.Range("H132").Value = 4663.5
$ws.Range("I132").Value = 4663.5
$ws.Range("K132").Value = 13990.5
$ws.Range("M132").Value = -11460.5
$ws.Range("H136").Value = 1353.6957
$ws.Range("I136").Value = 1292.1428
$ws.Range("K136").Value = 3876.4284
$ws.Range("M136").Value = -1326.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 461.66666
$ws.Range("I22").Value = 461.66666
$ws.Range("K22").Value = 461.66666
$ws.Range("M22").Value = -288.66666
$ws.Range("H92").Value = 36912.75
$ws.Range("J92").Value = 36912.75
$ws.Range("L92").Value = 36912.75
$ws.Range("N92").Value = -41904.75
$ws.Range("H95").Value = 6380.75
$ws.Range("J95").Value = 6380.75
$ws.Range("L95").Value = 6380.75
$ws.Range("N95").Value = -11872.75
$ws.Range("H105").Value = 4061.524
$ws.Range("I105").Value = 3215.7273
$ws.Range("K105").Value = 3215.7273
$ws.Range("M105").Value = -1468.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 7450
$ws.Range("I69").Value = 7450
$ws.Range("K69").Value = 7450
$ws.Range("M69").Value = -6701
$ws.Range("H72").Value = 7450
$ws.Range("I72").Value = 7450
$ws.Range("K72").Value = 22350
$ws.Range("M72").Value = -18606
$ws.Range("H122").Value = 4999.2
$ws.Range("I122").Value = 5249
$ws.Range("K122").Value = 15747
$ws.Range("M122").Value = -13297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 90.31579000000001
$ws.Range("J2").Value = 103.375
$ws.Range("L2").Value = 620.25
$ws.Range("N2").Value = -846.25
$ws.Range("H17").Value = 221
$ws.Range("I17").Value = 77.625
$ws.Range("K17").Value = 232.875
$ws.Range("M17").Value = -63.875
$ws.Range("H23").Value = 724.6667
$ws.Range("J23").Value = 694
$ws.Range("L23").Value = 2082
$ws.Range("N23").Value = -2552
$ws.Range("H26").Value = 1347.4814
$ws.Range("I26").Value = 1512.6086
$ws.Range("J26").Value = 398
$ws.Range("K26").Value = 4537.825800000001
$ws.Range("L26").Value = 1194
$ws.Range("M26").Value = -4249.825800000001
$ws.Range("N26").Value = -1770
$ws.Range("H68").Value = 541.5
$ws.Range("I68").Value = 497.33334
$ws.Range("J68").Value = 674
$ws.Range("K68").Value = 1492.00002
$ws.Range("L68").Value = 2022
$ws.Range("M68").Value = -681.0000199999999
$ws.Range("N68").Value = -3644
$ws.Range("H69").Value = 4508.8237
$ws.Range("I69").Value = 930
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 2790
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -1979
$ws.Range("N69").Value = -19622
$ws.Range("H71").Value = 541.5
$ws.Range("I71").Value = 497.33334
$ws.Range("J71").Value = 674
$ws.Range("K71").Value = 4476.00006
$ws.Range("L71").Value = 6066
$ws.Range("M71").Value = -420.0000600000003
$ws.Range("N71").Value = -14178
$ws.Range("H72").Value = 4508.8237
$ws.Range("I72").Value = 930
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 8370
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -4314
$ws.Range("N72").Value = -62112
$ws.Range("H92").Value = 278.2857
$ws.Range("J92").Value = 310.25
$ws.Range("L92").Value = 930.75
$ws.Range("N92").Value = -3426.75
$ws.Range("H113").Value = 1608.4
$ws.Range("J113").Value = 1561.3334
$ws.Range("L113").Value = 4684.0002
$ws.Range("N113").Value = -9024.0002
$ws.Range("H122").Value = 545.125
$ws.Range("I122").Value = 577
$ws.Range("J122").Value = 492
$ws.Range("K122").Value = 5193
$ws.Range("L122").Value = 4428
$ws.Range("M122").Value = -2743
$ws.Range("N122").Value = -9328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 38996.668
$ws.Range("J106").Value = 38996.668
$ws.Range("L106").Value = 38996.668
$ws.Range("N106").Value = -41520.668
$ws.Range("H122").Value = 4999.3335
$ws.Range("I122").Value = 4999
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14997
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12547
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39990
$ws.Range("I15").Value = 39990
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 39990
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -39702
$ws.Range("H39").Value = 25000
$ws.Range("I39").Value = 25000
$ws.Range("K39").Value = 25000
$ws.Range("M39").Value = -24587
$ws.Range("H41").Value = 16994.5
$ws.Range("J41").Value = 14659.667
$ws.Range("L41").Value = 14659.667
$ws.Range("N41").Value = -15439.667
$ws.Range("H45").Value = 45969.168
$ws.Range("I45").Value = 29967
$ws.Range("J45").Value = 53970.25
$ws.Range("K45").Value = 29967
$ws.Range("L45").Value = 53970.25
$ws.Range("M45").Value = -29476
$ws.Range("N45").Value = -54952.25
$ws.Range("H96").Value = 1499.5
$ws.Range("J96").Value = 1500
$ws.Range("L96").Value = 1500
$ws.Range("N96").Value = -4246
$ws.Range("H133").Value = 42499.5
$ws.Range("J133").Value = 42499.5
$ws.Range("L133").Value = 42499.5
$ws.Range("N133").Value = -52619.5
$ws.Range("N15").ClearContents()
